# Update cryptos list figures (price & 1h volume change) and reorder two
# coin pairs (VeChain/HuobiToken and FraxShare/InjectiveProtocol swap rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.768.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +6.23%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.053.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.64%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'252.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +4.50%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.653"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.54%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'65.54"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +16.28%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +5.94%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'59.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.57%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.81%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.09%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.928"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.92%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'15.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +8.07%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.351.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.53%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +6.92%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'20.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +20.42%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.046.49"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'37.635.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +6.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'73.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +5.33%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0881"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +5.75%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +6.34%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'238.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.42%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +15.67%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.04%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +5.17%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +6.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'160.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.93%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'19.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +2.68%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.66%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.97%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +24.95%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +7.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.76"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +12.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0617"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.64%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +4.46%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +3.42%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.08%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +21.91%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +15.83%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +24.94%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +4.84%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'VeChain"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.0220"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +5.44%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'HuobiToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.07%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +5.87%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'FraxShare"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'8.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +10.27%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'17.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +10.52%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'95.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +5.78%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.429.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.06%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'47.63"
$ws.Range("D51").Style = "Normal"
